$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Mark the pseudoinverse waterfall-chart picture run as <w:noProof/>.
#    This is "Picture 8" (2994660 x 2982024 EMU == 235.8 x 234.8 pt).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    if ([Math]::Round($shp.Width, 1) -eq 235.8 -and [Math]::Round($shp.Height, 1) -eq 234.8) {
        $shp.Range.NoProofing = 1
    }
}

# ---------------------------------------------------------------------------
# 2. Insert the new "We are providing two solutions..." / "Solution 1 ..."
#    paragraphs right before the "If we have A as 3x2 ..." paragraph (which
#    immediately follows the "V - Right Singular Vector , square matrix"
#    paragraph).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Right Singular Vector , square matrix", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$vPara = $rng.Paragraphs(1)
$targetPara = $vPara.Next()

$insertPoint = $targetPara.Range
$insertPoint.Collapse(1)
$insertPoint.InsertParagraphBefore()
$insertPoint.InsertParagraphBefore()
$insertPoint.InsertParagraphBefore()
$insertPoint.InsertParagraphBefore()

$p1 = $vPara.Next()
$p2 = $p1.Next()
$p3 = $p2.Next()
$p4 = $p3.Next()

$p2.Range.Text = "We are providing two solutions – one using a random 3x2 matrix and other using a housing dataset"

$p3r = $p3.Range
$p3r.Collapse(1)
$p3r.InsertAfter("Solution 1 – random 3x2 matrix ")
$p3.Range.InsertAfter(". ")

# ---------------------------------------------------------------------------
# 3. "If we have A as 3x2 , " -> " we have A as 3x2 , " (drop the leading "If").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("If we have A as 3x2 , ", $true, $false, $false, $false, $false, $true, 1, $false, " we have A as 3x2 , ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Add "Solution 2 – using Housing dataset" to the empty paragraph that
#    follows "V shape will be 2x2".
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("V shape will be 2x2", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$vShapePara = $rng2.Paragraphs(1)
$sol2Para = $vShapePara.Next()
$sol2Para.Range.Text = "Solution 2 – using Housing dataset"
